$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.08153196399554162
$ws.Range("B3").Value = -0.002531718893118817
$ws.Range("C3").Value = 0.0005331245300238389
$ws.Range("D3").Value = -3.779610287713201
$ws.Range("E3").Value = 0.02395263000825441
$ws.Range("F3").Value = -0.003576626419837724
$ws.Range("G3").Value = -0.001486811366399911
$ws.Range("H3").Value = 0.0790002451024228
$ws.Range("B4").Value = -0.00243516539028261
$ws.Range("C4").Value = 0.0007689672475205888
$ws.Range("D4").Value = -2.471102375843091
$ws.Range("E4").Value = 0.1720778460165632
$ws.Range("F4").Value = -0.003942317176089138
$ws.Range("G4").Value = -0.0009280136044760815
$ws.Range("H4").Value = 0.07909679860525901
$ws.Range("B5").Value = 0.01639614775064075
$ws.Range("C5").Value = 0.001670676520155925
$ws.Range("D5").Value = 7.209847351451951
$ws.Range("E5").Value = 0.1050469992975436
$ws.Range("F5").Value = 0.01312167348506571
$ws.Range("G5").Value = 0.01967062201621579
$ws.Range("H5").Value = 0.09792811174618238
$ws.Range("B6").Value = 0.02606912521302118
$ws.Range("C6").Value = 0.002397724710856256
$ws.Range("D6").Value = 16.34108086660036
$ws.Range("E6").Value = 0.06482583496268202
$ws.Range("F6").Value = 0.02136966179019903
$ws.Range("G6").Value = 0.03076858863584334
$ws.Range("H6").Value = 0.1076010892085628
$ws.Range("B7").Value = 0.06194933849454791
$ws.Range("C7").Value = 0.002702056733527469
$ws.Range("D7").Value = 19.02017140743589
$ws.Range("E7").Value = 0.004183796074436664
$ws.Range("F7").Value = 0.0566533943456923
$ws.Range("G7").Value = 0.06724528264340351
$ws.Range("H7").Value = 0.1434813024900895
$ws.Range("B8").Value = 0.07832023662514136
$ws.Range("C8").Value = 0.004773292697065111
$ws.Range("D8").Value = 15.22462831511043
$ws.Range("E8").Value = [double]"2.167649901899604e-05"
$ws.Range("F8").Value = 0.06896473453858969
$ws.Range("G8").Value = 0.08767573871169303
$ws.Range("H8").Value = 0.159852200620683
$ws.Range("B9").Value = 0.08952563088236325
$ws.Range("C9").Value = 0.00361650455291395
$ws.Range("D9").Value = 19.11372881127244
$ws.Range("E9").Value = [double]"1.388182492630752e-13"
$ws.Range("F9").Value = 0.08243739716695214
$ws.Range("G9").Value = 0.09661386459777435
$ws.Range("H9").Value = 0.1710575948779049
$ws.Range("B10").Value = -0.08153196399554162
$ws.Range("C10").Value = 0.0003939230741259212
$ws.Range("D10").Value = -218.0836416302556
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = -0.08230404106765341
$ws.Range("G10").Value = -0.08075988692342985
$ws.Range("B11").Value = -0.02763911866766634
$ws.Range("C11").Value = 0.0004546637769078559
$ws.Range("D11").Value = -61.38753894431049
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = -0.02853024568659471
$ws.Range("G11").Value = -0.02674799164873797
$ws.Range("H11").Value = 0.05389284532787528
$ws.Range("B12").Value = -0.0255228205810788
$ws.Range("C12").Value = 0.00044319023192568
$ws.Range("D12").Value = -57.87389626801608
$ws.Range("E12").Value = [double]"4.061180681985198e-252"
$ws.Range("F12").Value = -0.02639145980376092
$ws.Range("G12").Value = -0.02465418135839667
$ws.Range("H12").Value = 0.05600914341446282
$ws.Range("B13").Value = -0.01844018663834551
$ws.Range("C13").Value = 0.0004402634228673703
$ws.Range("D13").Value = -41.824774319745
$ws.Range("E13").Value = [double]"5.830736950348806e-137"
$ws.Range("F13").Value = -0.01930308940757593
$ws.Range("G13").Value = -0.0175772838691151
$ws.Range("H13").Value = 0.06309177735719611
$ws.Range("B14").Value = -0.01486566938585221
$ws.Range("C14").Value = 0.0004282648876418308
$ws.Range("D14").Value = -34.38204691329545
$ws.Range("E14").Value = [double]"4.103151666307971e-94"
$ws.Range("F14").Value = -0.01570505537331114
$ws.Range("G14").Value = -0.01402628339839328
$ws.Range("H14").Value = 0.06666629460968941
$ws.Range("B15").Value = -0.01184864904336384
$ws.Range("C15").Value = 0.0004105191210916494
$ws.Range("D15").Value = -28.91165082432362
$ws.Range("E15").Value = [double]"3.733402186643536e-60"
$ws.Range("F15").Value = -0.01265325383808668
$ws.Range("G15").Value = -0.01104404424864101
$ws.Range("H15").Value = 0.06968331495217778
$ws.Range("B16").Value = -0.008174794833566862
$ws.Range("C16").Value = 0.0003944960241504752
$ws.Range("D16").Value = -21.20696961925823
$ws.Range("E16").Value = [double]"6.160787758623587e-16"
$ws.Range("F16").Value = -0.008947994852918447
$ws.Range("G16").Value = -0.007401594814215281
$ws.Range("H16").Value = 0.07335716916197475
$ws.Range("B17").Value = -0.006881754817985328
$ws.Range("C17").Value = 0.0004122856760108743
$ws.Range("D17").Value = -17.06414040222145
$ws.Range("E17").Value = [double]"6.911515476269726e-05"
$ws.Range("F17").Value = -0.007689822041765011
$ws.Range("G17").Value = -0.006073687594205646
$ws.Range("H17").Value = 0.0746502091775563
$ws.Range("B18").Value = -0.005589551231773072
$ws.Range("C18").Value = 0.0004273191206168012
$ws.Range("D18").Value = -13.12663135059598
$ws.Range("E18").Value = [double]"3.642439063838546e-11"
$ws.Range("F18").Value = -0.006427083562127928
$ws.Range("G18").Value = -0.004752018901418216
$ws.Range("H18").Value = 0.07594241276376855
$ws.Range("B19").Value = -0.00481390789803541
$ws.Range("C19").Value = 0.0004238454221650106
$ws.Range("D19").Value = -11.14499012116771
$ws.Range("E19").Value = 0.002230726533510322
$ws.Range("F19").Value = -0.005644631880210525
$ws.Range("G19").Value = -0.003983183915860294
$ws.Range("H19").Value = 0.0767180560975062
$ws.Range("B20").Value = -0.003825348524689179
$ws.Range("C20").Value = 0.0004217804948601853
$ws.Range("D20").Value = -9.405915846298001
$ws.Range("E20").Value = 0.0004993686053760845
$ws.Range("F20").Value = -0.004652025282463954
$ws.Range("G20").Value = -0.002998671766914403
$ws.Range("H20").Value = 0.07770661547085245
$ws.Range("B21").Value = -0.003205876766728889
$ws.Range("C21").Value = 0.0004245861671905199
$ws.Range("D21").Value = -7.49696308564609
$ws.Range("E21").Value = 0.03684562912507663
$ws.Range("F21").Value = -0.004038052546305702
$ws.Range("G21").Value = -0.002373700987152077
$ws.Range("H21").Value = 0.07832608722881274
$ws.Range("B22").Value = -0.003167623101848114
$ws.Range("C22").Value = 0.0004199117778642859
$ws.Range("D22").Value = -7.862341003561978
$ws.Range("E22").Value = 0.1202869288539428
$ws.Range("F22").Value = -0.003990637201269249
$ws.Range("G22").Value = -0.00234460900242698
$ws.Range("H22").Value = 0.07836434089369351
$ws.Range("B23").Value = -0.002366041515784302
$ws.Range("C23").Value = 0.0004180229709331765
$ws.Range("D23").Value = -5.972833798911966
$ws.Range("E23").Value = 0.1953375207126779
$ws.Range("F23").Value = -0.003185353587954022
$ws.Range("G23").Value = -0.001546729443614581
$ws.Range("H23").Value = 0.07916592247975732
$ws.Range("B24").Value = -0.001475877131171742
$ws.Range("C24").Value = 0.0004052090406180862
$ws.Range("D24").Value = -4.960482890023501
$ws.Range("E24").Value = 0.2056482850097688
$ws.Range("F24").Value = -0.002270074269457065
$ws.Range("G24").Value = -0.0006816799928864178
$ws.Range("H24").Value = 0.08005608686436988
$ws.Range("B25").Value = -0.002016819415640612
$ws.Range("C25").Value = 0.0003883084376025476
$ws.Range("D25").Value = -6.697891784284449
$ws.Range("E25").Value = 0.001453500013514659
$ws.Range("F25").Value = -0.002777891870949813
$ws.Range("G25").Value = -0.001255746960331411
$ws.Range("H25").Value = 0.07951514457990101
$ws.Range("B26").Value = 0.135366388148127
$ws.Range("C26").Value = 0.002573534430795923
$ws.Range("D26").Value = 42.44227064394657
$ws.Range("E26").Value = [double]"1.67010038170914e-43"
$ws.Range("F26").Value = 0.1303223420097282
$ws.Range("G26").Value = 0.1404104342865257
$ws.Range("H26").Value = 0.2168983521436686
